# Apply the target edits to the "Gasto Funcionamiento" sheet:
#  - Row 12 (Materiales y suministros detail): change item description + total (E12)
#  - Row 13 (was the ABRAZADERA detail line): replace with a bare "Servicios" header row
#  - Old row 14 (bare "Servicios" header) and row 15 (MANTENIMIENTO DE VEHICULOS GENERALES
#    detail line) are removed outright, per the commit note that area maintenance detail
#    is pending.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gasto Funcionamiento")

# Row 12: swap the line-item description and correct the total (column E).
$ws.Range("B12").Value = "MACHO P/R TUBO NPT DE 1.1/2  X 11.5 HILOS"
$ws.Range("E12").Value = 655.83

# Row 13 currently holds the "ABRAZADERA ..." detail line (B..R populated). Clear it out
# and delete rows 14-15 (the old "Servicios" header + its maintenance detail line), then
# put a bare "Servicios" header back on row 13 - mirroring row 11's bare "Materiales y
# suministros" header - leaving detail rows for that section to be filled in later.
$ws.Range("C13:R13").ClearContents()
$ws.Rows("14:15").Delete()
$ws.Range("B13").Value = "Servicios"
